# Update "想去人数" (number of people interested) values in the
# "展览" (rId1 / sheet1) and "全部类型" (rId4 / sheet4) worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 759
$ws1.Range("F4").Value = 52
$ws1.Range("F5").Value = 28
$ws1.Range("F6").Value = 257
$ws1.Range("F7").Value = 3516
$ws1.Range("F8").Value = 72
$ws1.Range("F9").Value = 4157
$ws1.Range("F11").Value = 1042
$ws1.Range("F12").Value = 48

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 759
$ws4.Range("F4").Value = 52
$ws4.Range("F5").Value = 28
$ws4.Range("F7").Value = 257
$ws4.Range("F8").Value = 3516
$ws4.Range("F9").Value = 72
$ws4.Range("F10").Value = 4157
$ws4.Range("F12").Value = 1042
$ws4.Range("F13").Value = 48
